$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the issue date/time stamp in the general information section (A5)
$ws.Range("A5").Value = "Issue date: 02/12/2020 12:20:32"

# Update the Doc2VecTransfomer accuracy scores (row 17: MLP, SVC, LR, RF)
$ws.Range("G17").Value = "65.03*"
$ws.Range("H17").Value = "64.85*"
$ws.Range("I17").Value = "63.25*"
$ws.Range("J17").Value = "64.8*"

# Swap the highlight colors between the MLP (G17) and SVC (H17) result cells:
# G17 goes from blue to red, H17 goes from red to blue.
$ws.Range("G17").Font.Color = 255
$ws.Range("H17").Font.Color = 16711680
